$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-96 down to 12-97
$ws.Rows(11).Insert()

# Populate the new row 11 with the new record
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(11, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(11, 3).Value = "Bíobío"
$ws.Cells.Item(11, 4).Value = 44602
$ws.Cells.Item(11, 4).Style = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 100112032
$ws.Cells.Item(11, 7).Value = "Zapallo italiano"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 170
$ws.Cells.Item(11, 11).Value = 7000
$ws.Cells.Item(11, 12).Value = 7500
$ws.Cells.Item(11, 13).Value = 7235
$ws.Cells.Item(11, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 145
$ws.Cells.Item(11, 17).Value = 50
$ws.Cells.Item(11, 18).Value = "Hortaliza"
